# Faculty vaccination tracker header row: FSN | Name | Age | Phone | Vaccine_Dose
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "FSN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen the Vaccine_Dose column (target stored width ~14.18 chars; engine snaps
# column widths to an internal pixel grid, so 13.3 is the input that lands on
# the closest reachable stored width).
$ws.Columns.Item(5).ColumnWidth = 13.3

$ws.Range("F4").Select() | Out-Null
